# The deck's single slide master ("theme2.xml" in the package) currently uses
# the "Integral" theme color scheme; the notes master ("theme1.xml") currently
# uses the default "Office Theme" color scheme. The edit swaps these two
# color schemes so the slide master becomes "Office Theme" colors and the
# notes master becomes "Integral" colors.
#
# PowerPoint's object model doesn't expose a raw "swap these two theme parts"
# call, so we reproduce the effect by writing the target RGB values straight
# into the live theme's 12-slot ThemeColorScheme (dk1/lt1/dk2/lt2/accent1-6/
# hlink/folHlink) via $master.Theme.ThemeColorScheme.Colors(n).RGB -- this is
# the only COM surface that persists back into the OOXML theme part without
# clobbering the scheme's name attribute.

function ToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0,2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2,2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4,2), 16)
    # VBA/COM RGB() packs components as 0xBBGGRR.
    return ($b * 65536) + ($g * 256) + $r
}

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

# Target = the stock "Office Theme" color scheme (12 slots, in Colors() order).
$officeTheme = @(
    "000000",  # 1  dk1      - Text 1
    "FFFFFF",  # 2  lt1      - Background 1
    "44546A",  # 3  dk2      - Text 2
    "E7E6E6",  # 4  lt2      - Background 2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $officeTheme.Count; $i++) {
    $themeColors.Colors($i).RGB = ToRgbInt($officeTheme[$i - 1])
}
